$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.795.09"
$ws.Range("E2").Value = "  +4.93%  "
$ws.Range("D3").Value = "3.111.44"
$ws.Range("E3").Value = "  +3.03%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.62%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.103.65"
$ws.Range("E8").Value = "  +3.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.149"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.22%  "
$ws.Range("E11").Value = "  +7.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.467"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.45%  "
$ws.Range("E13").Value = "  +6.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.70%  "
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").Value = "3.630.24"
$ws.Range("E16").Value = "  +3.02%  "
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D18").Value = "3.113.88"
$ws.Range("E18").Value = "  +3.11%  "
$ws.Range("D19").Value = "62.801.47"
$ws.Range("E19").Value = "  +4.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "464.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.728"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.81%  "
$ws.Range("E23").Value = "  +5.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("E28").Value = "  +4.86%  "
$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.28%  "
$ws.Range("E33").Value = "  +7.43%  "
$ws.Range("E34").Value = "  +7.34%  "
$ws.Range("E35").Value = "  +11.21%  "
$ws.Range("E36").Value = "  +3.36%  "
$ws.Range("E37").Value = "  +2.02%  "
$ws.Range("E38").Value = "  +15.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "51.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "431.98"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.22%  "
$ws.Range("D42").Value = "2.938.51"
$ws.Range("E42").Value = "  +5.45%  "
$ws.Range("E43").Value = "  +4.19%  "
$ws.Range("E44").Value = "  +9.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.110"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.47%  "
